$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: username/password changed
$ws.Range("A3").Value = "test1234@gmail.com"
$ws.Range("B3").Value = "1234"
$ws.Range("C3").Value = "Valid"

# Row 6: username changed
$ws.Range("A6").Value = "test23456@gmail.com"

# Row 7: clear B7 and C7 values but keep formatting
$ws.Range("B7").ClearContents()
$ws.Range("C7").ClearContents()

# Update selection to A6
$ws.Range("A6").Select()
